$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update country names (shared string reorder) and timestamp
$ws.Range("A38").Value = "Indonesia"
$ws.Range("A39").Value = "Arabia Saudita"
$ws.Range("A40").Value = "Grecia"
$ws.Range("A41").Value = "Crucero"
$ws.Range("A42").Value = "Sudafrica"
$ws.Range("A43").Value = "Rusia"
$ws.Range("A44").Value = "Islandia"
$ws.Range("A45").Value = "Filipinas"
$ws.Range("A46").Value = "India"
$ws.Range("A47").Value = "Singapur"
$ws.Range("A48").Value = "Catar"
$ws.Range("A54").Value = "Estonia"
$ws.Range("A55").Value = "Egipto"
$ws.Range("A56").Value = "Barein"
$ws.Range("A57").Value = "Argentina"
$ws.Range("A58").Value = "Hong Kong"
$ws.Range("A59").Value = "Colombia"
$ws.Range("A72").Value = "Eslovaquia"
$ws.Range("A73").Value = "Nueva Zelanda"
$ws.Range("A123").Value = "Honduras"
$ws.Range("A124").Value = "Mayotte"
$ws.Range("A128").Value = "Polinesia Francesa"
$ws.Range("A129").Value = "Kenia"
$ws.Range("A143").Value = "Etiopia"
$ws.Range("A144").Value = "Tanzania"
$ws.Range("A151").Value = "Seychelles"
$ws.Range("A153").Value = "Dominica"
$ws.Range("A166").Value = "Congo"
$ws.Range("A167").Value = "Guinea"
$ws.Range("A169").Value = "Zambia"
$ws.Range("A170").Value = "Angola"
$ws.Range("A171").Value = "Mozambique"
$ws.Range("A172").Value = "Birmania"
$ws.Range("A173").Value = "Liberia"
$ws.Range("A174").Value = "Antigua y Barbuda"
$ws.Range("A175").Value = "Santa Lucia"
$ws.Range("A176").Value = "San Bartolome"
$ws.Range("A177").Value = "Republica del Chad"
$ws.Range("A178").Value = "Republica de Yibuti"
$ws.Range("A179").Value = "Republica de Africa Central"
$ws.Range("A180").Value = "Gambia"
$ws.Range("A181").Value = "Nepal"
$ws.Range("A183").Value = "Sudan"
$ws.Range("A184").Value = "Cabo Verde"
$ws.Range("A185").Value = "Mauritania"
$ws.Range("A186").Value = "San Martin (Parte Holandesa)"
$ws.Range("A190").Value = "San Vicente y las Granadinas"
$ws.Range("A191").Value = "Granada"
$ws.Range("A192").Value = "Siria"
$ws.Range("A193").Value = "Montserrat"
$ws.Range("A194").Value = "Libia"
$ws.Range("A195").Value = "Eritrea"
$ws.Range("A196").Value = "Timor Oriental"
$ws.Range("A197").Value = "Islas Turcas y Caicos"
$ws.Range("A198").Value = "Belice"
$ws.Range("A200").Value = "Papua Nueva Guinea"
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 10:16"

# Update numeric stats for rows whose data changed
$ws.Range("B15").Value = 5448
$ws.Range("C15").Value = 165
$ws.Range("E15").Value = 5409
$ws.Range("B21").Value = 2300
$ws.Range("C21").Value = 1
$ws.Range("E21").Value = 2244
$ws.Range("F21").Value = 144
$ws.Range("B26").Value = 1715
$ws.Range("C26").Value = 124
$ws.Range("E26").Value = 1682
$ws.Range("B38").Value = 790
$ws.Range("C38").Value = 104
$ws.Range("D38").Value = 31
$ws.Range("E38").Value = 701
$ws.Range("G38").Value = 3
$ws.Range("H38").Value = 58
$ws.Range("B39").Value = 767
$ws.Range("D39").Value = 28
$ws.Range("E39").Value = 738
$ws.Range("F39").Value = 0
$ws.Range("H39").Value = 1
$ws.Range("B40").Value = 743
$ws.Range("D40").Value = 29
$ws.Range("E40").Value = 694
$ws.Range("F40").Value = 35
$ws.Range("H40").Value = 20
$ws.Range("B41").Value = 712
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 587
$ws.Range("E41").Value = 115
$ws.Range("F41").Value = 15
$ws.Range("H41").Value = 10
$ws.Range("B42").Value = 709
$ws.Range("C42").Value = 155
$ws.Range("D42").Value = 4
$ws.Range("E42").Value = 705
$ws.Range("F42").Value = 2
$ws.Range("H42").Value = 0
$ws.Range("B43").Value = 658
$ws.Range("C43").Value = 163
$ws.Range("D43").Value = 29
$ws.Range("E43").Value = 628
$ws.Range("F43").Value = 8
$ws.Range("H43").Value = 1
$ws.Range("B44").Value = 648
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 51
$ws.Range("E44").Value = 595
$ws.Range("F44").Value = 13
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 2
$ws.Range("B45").Value = 636
$ws.Range("C45").Value = 84
$ws.Range("D45").Value = 26
$ws.Range("E45").Value = 572
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 3
$ws.Range("H45").Value = 38
$ws.Range("B46").Value = 562
$ws.Range("C46").Value = 26
$ws.Range("D46").Value = 40
$ws.Range("E46").Value = 512
$ws.Range("F46").Value = 0
$ws.Range("H46").Value = 10
$ws.Range("B47").Value = 558
$ws.Range("D47").Value = 156
$ws.Range("E47").Value = 400
$ws.Range("F47").Value = 17
$ws.Range("H47").Value = 2
$ws.Range("B48").Value = 526
$ws.Range("D48").Value = 41
$ws.Range("E48").Value = 485
$ws.Range("F48").Value = 6
$ws.Range("H48").Value = 0
$ws.Range("B54").Value = 404
$ws.Range("C54").Value = 35
$ws.Range("D54").Value = 8
$ws.Range("E54").Value = 396
$ws.Range("F54").Value = 5
$ws.Range("H54").Value = 0
$ws.Range("B55").Value = 402
$ws.Range("D55").Value = 80
$ws.Range("E55").Value = 302
$ws.Range("F55").Value = 0
$ws.Range("H55").Value = 20
$ws.Range("B56").Value = 392
$ws.Range("D56").Value = 177
$ws.Range("E56").Value = 212
$ws.Range("F56").Value = 2
$ws.Range("H56").Value = 3
$ws.Range("D57").Value = 52
$ws.Range("E57").Value = 329
$ws.Range("F57").Value = 0
$ws.Range("H57").Value = 6
$ws.Range("B58").Value = 387
$ws.Range("D58").Value = 102
$ws.Range("E58").Value = 281
$ws.Range("F58").Value = 4
$ws.Range("H58").Value = 4
$ws.Range("B59").Value = 378
$ws.Range("D59").Value = 6
$ws.Range("E59").Value = 369
$ws.Range("F59").Value = 0
$ws.Range("H59").Value = 3
$ws.Range("E66").Value = 251
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 3
$ws.Range("B72").Value = 216
$ws.Range("C72").Value = 12
$ws.Range("D72").Value = 7
$ws.Range("E72").Value = 209
$ws.Range("F72").Value = 2
$ws.Range("B73").Value = 205
$ws.Range("D73").Value = 22
$ws.Range("E73").Value = 183
$ws.Range("F73").Value = 0
$ws.Range("F85").Value = 20
$ws.Range("E90").Value = 108
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 4
$ws.Range("D100").Value = 29
$ws.Range("E100").Value = 57
$ws.Range("C123").Value = 6
$ws.Range("C124").Value = 0
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 0
$ws.Range("H180").Value = 1
$ws.Range("C181").Value = 1
$ws.Range("D181").Value = 1
$ws.Range("H181").Value = 0
